$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efnb2"
$ws.Cells.Item(2,3).Value = "Ephb3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 36.899643
$ws.Cells.Item(2,8).Value = 110.698929
$ws.Cells.Item(2,9).Value = 0.7238945645409351
$ws.Cells.Item(2,10).Value = 0.7238945645409351
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.07588399999999999
$ws.Cells.Item(2,14).Value = 0.227652
$ws.Cells.Item(2,15).Value = 0.005197207581907009
$ws.Cells.Item(2,16).Value = 0.00519720758190701
$ws.Cells.Item(2,17).Value = 2.800092509411999
$ws.Cells.Item(2,18).Value = 25.200832584708
$ws.Cells.Item(2,19).Value = 0.003762230319333421
$ws.Cells.Item(2,20).Value = 0.003762230319333422

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efnb2"
$ws.Cells.Item(3,3).Value = "Ephb3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 36.899643
$ws.Cells.Item(3,8).Value = 110.698929
$ws.Cells.Item(3,9).Value = 0.7238945645409351
$ws.Cells.Item(3,10).Value = 0.7238945645409351
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 11.71993
$ws.Cells.Item(3,14).Value = 35.15979
$ws.Cells.Item(3,15).Value = 0.80268447967186
$ws.Cells.Item(3,16).Value = 0.80268447967186
$ws.Cells.Item(3,17).Value = 432.46123298499
$ws.Cells.Item(3,18).Value = 3892.15109686491
$ws.Cells.Item(3,19).Value = 0.5810589318758281
$ws.Cells.Item(3,20).Value = 0.5810589318758281

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efnb2"
$ws.Cells.Item(4,3).Value = "Ephb3"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 36.899643
$ws.Cells.Item(4,8).Value = 110.698929
$ws.Cells.Item(4,9).Value = 0.7238945645409351
$ws.Cells.Item(4,10).Value = 0.7238945645409351
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.805103666666666
$ws.Cells.Item(4,14).Value = 8.415310999999999
$ws.Cells.Item(4,15).Value = 0.1921183127462331
$ws.Cells.Item(4,16).Value = 0.1921183127462331
$ws.Cells.Item(4,17).Value = 103.507323877991
$ws.Cells.Item(4,18).Value = 931.5659149019189
$ws.Cells.Item(4,19).Value = 0.1390734023457736
$ws.Cells.Item(4,20).Value = 0.1390734023457736

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Efnb2"
$ws.Cells.Item(5,3).Value = "Ephb3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.374819
$ws.Cells.Item(5,8).Value = 10.124457
$ws.Cells.Item(5,9).Value = 0.0662069584361419
$ws.Cells.Item(5,10).Value = 0.0662069584361419
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.07588399999999999
$ws.Cells.Item(5,14).Value = 0.227652
$ws.Cells.Item(5,15).Value = 0.005197207581907009
$ws.Cells.Item(5,16).Value = 0.00519720758190701
$ws.Cells.Item(5,17).Value = 0.256094764996
$ws.Cells.Item(5,18).Value = 2.304852884964
$ws.Cells.Item(5,19).Value = 0.0003440913063593189
$ws.Cells.Item(5,20).Value = 0.000344091306359319

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Efnb2"
$ws.Cells.Item(6,3).Value = "Ephb3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.374819
$ws.Cells.Item(6,8).Value = 10.124457
$ws.Cells.Item(6,9).Value = 0.0662069584361419
$ws.Cells.Item(6,10).Value = 0.0662069584361419
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 11.71993
$ws.Cells.Item(6,14).Value = 35.15979
$ws.Cells.Item(6,15).Value = 0.80268447967186
$ws.Cells.Item(6,16).Value = 0.80268447967186
$ws.Cells.Item(6,17).Value = 39.55264244267
$ws.Cells.Item(6,18).Value = 355.97378198403
$ws.Cells.Item(6,19).Value = 0.05314329798297102
$ws.Cells.Item(6,20).Value = 0.05314329798297102

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efnb2"
$ws.Cells.Item(7,3).Value = "Ephb3"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3.374819
$ws.Cells.Item(7,8).Value = 10.124457
$ws.Cells.Item(7,9).Value = 0.0662069584361419
$ws.Cells.Item(7,10).Value = 0.0662069584361419
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.805103666666666
$ws.Cells.Item(7,14).Value = 8.415310999999999
$ws.Cells.Item(7,15).Value = 0.1921183127462331
$ws.Cells.Item(7,16).Value = 0.1921183127462331
$ws.Cells.Item(7,17).Value = 9.466717151236331
$ws.Cells.Item(7,18).Value = 85.20045436112699
$ws.Cells.Item(7,19).Value = 0.01271956914681156
$ws.Cells.Item(7,20).Value = 0.01271956914681156

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Efnb2"
$ws.Cells.Item(8,3).Value = "Ephb3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 10.699319
$ws.Cells.Item(8,8).Value = 32.097957
$ws.Cells.Item(8,9).Value = 0.2098984770229228
$ws.Cells.Item(8,10).Value = 0.2098984770229228
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.07588399999999999
$ws.Cells.Item(8,14).Value = 0.227652
$ws.Cells.Item(8,15).Value = 0.005197207581907009
$ws.Cells.Item(8,16).Value = 0.00519720758190701
$ws.Cells.Item(8,17).Value = 0.811907122996
$ws.Cells.Item(8,18).Value = 7.307164106964
$ws.Cells.Item(8,19).Value = 0.001090885956214269
$ws.Cells.Item(8,20).Value = 0.001090885956214269

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Efnb2"
$ws.Cells.Item(9,3).Value = "Ephb3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 10.699319
$ws.Cells.Item(9,8).Value = 32.097957
$ws.Cells.Item(9,9).Value = 0.2098984770229228
$ws.Cells.Item(9,10).Value = 0.2098984770229228
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 11.71993
$ws.Cells.Item(9,14).Value = 35.15979
$ws.Cells.Item(9,15).Value = 0.80268447967186
$ws.Cells.Item(9,16).Value = 0.80268447967186
$ws.Cells.Item(9,17).Value = 125.39526972767
$ws.Cells.Item(9,18).Value = 1128.55742754903
$ws.Cells.Item(9,19).Value = 0.1684822498130607
$ws.Cells.Item(9,20).Value = 0.1684822498130607

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Efnb2"
$ws.Cells.Item(10,3).Value = "Ephb3"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 10.699319
$ws.Cells.Item(10,8).Value = 32.097957
$ws.Cells.Item(10,9).Value = 0.2098984770229228
$ws.Cells.Item(10,10).Value = 0.2098984770229228
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.805103666666666
$ws.Cells.Item(10,14).Value = 8.415310999999999
$ws.Cells.Item(10,15).Value = 0.1921183127462331
$ws.Cells.Item(10,16).Value = 0.1921183127462331
$ws.Cells.Item(10,17).Value = 30.01269895773633
$ws.Cells.Item(10,18).Value = 270.114290619627
$ws.Cells.Item(10,19).Value = 0.04032534125364791
$ws.Cells.Item(10,20).Value = 0.04032534125364791
